$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Main Functions and Events")
$ws2 = $wb.Worksheets.Item("Changes to variables")

# -----------------------------------------------------------------
# 1) "Changes to variables": [m] balancesBonus -> ...[m] balancesBonus
# -----------------------------------------------------------------
$ws2.Range("B16").Value = "...[m] balancesBonus"

# -----------------------------------------------------------------
# 2) "Main Functions and Events": add a new function block
#    (addToWhitelist() / addToWhitelistParams()) above the
#    "cancelPending" block, pushing rows 20-26 down to 21-27.
# -----------------------------------------------------------------
$ws1.Rows("20:20").Insert()

$ws1.Range("A18").Value = "addToWhitelist()"
$ws1.Range("A19").Value = "addToWhitelistParams()"

# B19 picks up the same highlight fill as B18 (the "function name" cell)
$ws1.Range("B18").Copy()
$ws1.Range("B19").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# -----------------------------------------------------------------
# 3) Add a new trailing worksheet named "Sheet1"
# -----------------------------------------------------------------
$lastIndex = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($lastIndex)
$wb.Worksheets.Add($null, $lastSheet) | Out-Null

# -----------------------------------------------------------------
# 4) Restore selections / active sheet
# -----------------------------------------------------------------
$ws1.Select() | Out-Null
$ws1.Range("A19").Select() | Out-Null

$ws2.Select() | Out-Null
$ws2.Range("J7").Select() | Out-Null
